$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, forcing text storage so Excel
# does not reinterpret numeric-looking strings (e.g. "6.080", "1.000")
# as numbers and strip formatting/trailing zeros.
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextCell $ws.Range('D2') '29.035.54'
$ws.Range('E2').Value = '  +0.02%  '

# Row 3
Set-TextCell $ws.Range('D3') '1.831.32'
$ws.Range('E3').Value = '  +0.10%  '

# Row 4
Set-TextCell $ws.Range('D4') '0.9983'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
Set-TextCell $ws.Range('D5') '244.59'
$ws.Range('E5').Value = '  +1.47%  '

# Row 6
Set-TextCell $ws.Range('D6') '0.6322'
$ws.Range('E6').Value = '  +1.36%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
Set-TextCell $ws.Range('D8') '0.07516'
$ws.Range('E8').Value = '  -0.16%  '

# Row 9
$ws.Range('E9').Value = '  +1.04%  '

# Row 10
$ws.Range('E10').Value = '  +0.58%  '

# Row 11
Set-TextCell $ws.Range('D11') '0.07702'
$ws.Range('E11').Value = '  +0.78%  '

# Row 12
Set-TextCell $ws.Range('D12') '1.829.71'
$ws.Range('E12').Value = '  +0.03%  '

# Row 13
Set-TextCell $ws.Range('D13') '4.996'
$ws.Range('E13').Value = '  +0.81%  '

# Row 14
Set-TextCell $ws.Range('D14') '0.6704'
$ws.Range('E14').Value = '  +0.98%  '

# Row 15
$ws.Range('E15').Value = '  +0.89%  '

# Row 16
Set-TextCell $ws.Range('D16') '0.000009588'
$ws.Range('E16').Value = '  +5.82%  '

# Row 17
Set-TextCell $ws.Range('D17') '6.080'
$ws.Range('E17').Value = '  +1.26%  '

# Row 18
Set-TextCell $ws.Range('D18') '29.048.51'
$ws.Range('E18').Value = '  +0.39%  '

# Row 19
Set-TextCell $ws.Range('D19') '12.58'
$ws.Range('E19').Value = '  +2.05%  '

# Row 20
Set-TextCell $ws.Range('D20') '226.52'
$ws.Range('E20').Value = '  +0.76%  '

# Row 21
Set-TextCell $ws.Range('D21') '0.9987'
$ws.Range('E21').Value = '  -0.20%  '

# Row 22
Set-TextCell $ws.Range('D22') '7.159'
$ws.Range('E22').Value = '  -0.42%  '

# Row 23
$ws.Range('E23').Value = '  -0.11%  '

# Row 24
Set-TextCell $ws.Range('D24') '160.12'
$ws.Range('E24').Value = '  +0.52%  '

# Row 25
Set-TextCell $ws.Range('D25') '0.1408'
$ws.Range('E25').Value = '  +3.86%  '

# Row 26
$ws.Range('E26').Value = '  +1.87%  '

# Row 27
Set-TextCell $ws.Range('D27') '17.91'
$ws.Range('E27').Value = '  +0.49%  '

# Row 28
Set-TextCell $ws.Range('D28') '1.498'
$ws.Range('E28').Value = '  +0.46%  '

# Row 29
Set-TextCell $ws.Range('D29') '4.128'
$ws.Range('E29').Value = '  +1.91%  '

# Row 30
Set-TextCell $ws.Range('D30') '4.065'
$ws.Range('E30').Value = '  +0.79%  '

# Row 31
Set-TextCell $ws.Range('D31') '1.199'
$ws.Range('E31').Value = '  -0.05%  '

# Row 32
Set-TextCell $ws.Range('D32') '0.05384'
$ws.Range('E32').Value = '  +3.54%  '

# Row 33
Set-TextCell $ws.Range('D33') '1.859'
$ws.Range('E33').Value = '  +1.23%  '

# Row 34
Set-TextCell $ws.Range('D34') '0.7443'
$ws.Range('E34').Value = '  +1.71%  '

# Row 35
$ws.Range('E35').Value = '  -1.09%  '

# Row 36
Set-TextCell $ws.Range('D36') '2.655'
$ws.Range('E36').Value = '  +1.78%  '

# Row 37
Set-TextCell $ws.Range('D37') '1.243.88'
$ws.Range('E37').Value = '  -3.17%  '

# Row 38
Set-TextCell $ws.Range('D38') '2.766'
$ws.Range('E38').Value = '  +0.25%  '

# Row 39
$ws.Range('E39').Value = '  +0.42%  '

# Row 40
Set-TextCell $ws.Range('D40') '6.647'
$ws.Range('E40').Value = '  +4.14%  '

# Row 41
Set-TextCell $ws.Range('D41') '0.9024'
$ws.Range('E41').Value = '  +1.06%  '

# Row 42
Set-TextCell $ws.Range('D42') '1.000'
$ws.Range('E42').Value = '  -0.12%  '

# Row 43
Set-TextCell $ws.Range('D43') '102.04'
$ws.Range('E43').Value = '  +0.78%  '

# Row 44
Set-TextCell $ws.Range('D44') '1.977.01'
$ws.Range('E44').Value = '  -0.17%  '

# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range('D45') '64.84'
$ws.Range('E45').Value = '  +2.12%  '

# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws.Range('D46') '0.00000000121'
$ws.Range('E46').Value = '  +1.74%  '

# Row 47
Set-TextCell $ws.Range('D47') '0.5105'
$ws.Range('E47').Value = '  -0.14%  '

# Row 48
Set-TextCell $ws.Range('D48') '0.4072'
$ws.Range('E48').Value = '  +2.51%  '

# Row 49
Set-TextCell $ws.Range('D49') '8.960'
$ws.Range('E49').Value = '  +0.91%  '

# Row 50
$ws.Range('E50').Value = '  -0.25%  '

# Row 51
$ws.Range('E51').Value = '  +0.33%  '
